$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: new BB1 cell (date serial), copy header style (bold/border/date fmt) from BA1
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$ws.Range("BB1").Value = 45986

# Rows 2-82: new BB column, values (rows 2-70 mirror column BA; rows 71-82 are revised forecast values)
$ws.Range("BB2").Value = 0.4011221634494007
$ws.Range("BB3").Value = 0.4238702649943349
$ws.Range("BB4").Value = 1.065831052964711
$ws.Range("BB5").Value = 0.1679892790148898
$ws.Range("BB6").Value = 0.7826320556679605
$ws.Range("BB7").Value = 1.608891637697312
$ws.Range("BB8").Value = 0.7314738238097078
$ws.Range("BB9").Value = 0.6066829402737994
$ws.Range("BB10").Value = 0.3745978975504585
$ws.Range("BB11").Value = 1.488552458234295
$ws.Range("BB12").Value = -1.412768783584468
$ws.Range("BB13").Value = 1.150059426400716
$ws.Range("BB14").Value = 0.08121742761659334
$ws.Range("BB15").Value = -0.02717560123848273
$ws.Range("BB16").Value = 0.7243763014332956
$ws.Range("BB17").Value = 0.06292774364702325
$ws.Range("BB18").Value = 0.6288730841304044
$ws.Range("BB19").Value = -0.04932303278563666
$ws.Range("BB20").Value = 0.7008158264609534
$ws.Range("BB21").Value = -0.05880790937801805
$ws.Range("BB22").Value = 0.4315408898675201
$ws.Range("BB23").Value = 0.01220354597239748
$ws.Range("BB24").Value = -0.0390655612641666
$ws.Range("BB25").Value = 0.5861016707528819
$ws.Range("BB26").Value = -0.09710917920445183
$ws.Range("BB27").Value = 0.5608392729090639
$ws.Range("BB28").Value = 0.4327290387686133
$ws.Range("BB29").Value = 0.6415186033097626
$ws.Range("BB30").Value = 0.3424960270584307
$ws.Range("BB31").Value = 0.450007896627298
$ws.Range("BB32").Value = 0.5281644231130258
$ws.Range("BB33").Value = 0.9569320119344127
$ws.Range("BB34").Value = 1.458965617195048
$ws.Range("BB35").Value = 1.616667503385756
$ws.Range("BB36").Value = 0.7391298723574948
$ws.Range("BB37").Value = 0.05368498059650051
$ws.Range("BB38").Value = 0.2682953781150843
$ws.Range("BB39").Value = 0.07985748687143257
$ws.Range("BB40").Value = 0.5497494968945205
$ws.Range("BB41").Value = 0.5291048039492949
$ws.Range("BB42").Value = 0.4298198366907684
$ws.Range("BB43").Value = -0.00001256079109168695
$ws.Range("BB44").Value = 0.5955762646739942
$ws.Range("BB45").Value = 0.1202567733712812
$ws.Range("BB46").Value = 0.378819647862997
$ws.Range("BB47").Value = 1.1
$ws.Range("BB48").Value = 0.5
$ws.Range("BB49").Value = 1.4
$ws.Range("BB50").Value = 0.1
$ws.Range("BB51").Value = 1.14595180410879
$ws.Range("BB52").Value = 1.264351713671076
$ws.Range("BB53").Value = 0.9847891746257034
$ws.Range("BB54").Value = 0.1306165969516826
$ws.Range("BB55").Value = -0.702620495601451
$ws.Range("BB56").Value = 4.902896932015494
$ws.Range("BB57").Value = -3.137970849842105
$ws.Range("BB58").Value = 0.1542746451113572
$ws.Range("BB59").Value = 0.833545251617636
$ws.Range("BB60").Value = 0.3586840777896043
$ws.Range("BB61").Value = -1.11464790305466
$ws.Range("BB62").Value = 0.1937043270045251
$ws.Range("BB63").Value = -1.254351945331166
$ws.Range("BB64").Value = 0.05976999529040938
$ws.Range("BB65").Value = 1.365536114841987
$ws.Range("BB66").Value = 0.5893749959732304
$ws.Range("BB67").Value = 0.450398338429352
$ws.Range("BB68").Value = 1.478730722800918
$ws.Range("BB69").Value = 1.16573263774626
$ws.Range("BB70").Value = 0.4367731001324842
$ws.Range("BB71").Value = 0.1765865160815849
$ws.Range("BB72").Value = 0.2412052862208469
$ws.Range("BB73").Value = 0.768168485846715
$ws.Range("BB74").Value = 0.768168485846715
$ws.Range("BB75").Value = 0.768168485846715
$ws.Range("BB76").Value = 0.768168485846715
$ws.Range("BB77").Value = 0.768168485846715
$ws.Range("BB78").Value = 0.768168485846715
$ws.Range("BB79").Value = 0.768168485846715
$ws.Range("BB80").Value = 0.768168485846715
$ws.Range("BB81").Value = 0.768168485846715
$ws.Range("BB82").Value = 0.768168485846715

# New row 83: copy date-cell style from A82 into A83, then set values
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.768168485846715

